$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: Polonia
$ws.Cells.Item(30, 1).Value = "Polonia"
$ws.Cells.Item(30, 2).Value = 7918
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 866
$ws.Cells.Item(30, 5).Value = 6738
$ws.Cells.Item(30, 6).Value = 160
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 314

# Row 74: Armenia
$ws.Cells.Item(74, 1).Value = "Armenia"
$ws.Cells.Item(74, 2).Value = 1201
$ws.Cells.Item(74, 3).Value = 42
$ws.Cells.Item(74, 4).Value = 402
$ws.Cells.Item(74, 5).Value = 780
$ws.Cells.Item(74, 6).Value = 30
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 19

# Row 75: Bosnia y Herzegovina
$ws.Cells.Item(75, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(75, 2).Value = 1167
$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(75, 4).Value = 277
$ws.Cells.Item(75, 5).Value = 847
$ws.Cells.Item(75, 6).Value = 4
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 43

# Row 82: Afganistan
$ws.Cells.Item(82, 1).Value = "Afganistan"
$ws.Cells.Item(82, 2).Value = 906
$ws.Cells.Item(82, 3).Value = 66
$ws.Cells.Item(82, 4).Value = 99
$ws.Cells.Item(82, 5).Value = 777
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 30

# Row 83: Cuba
$ws.Cells.Item(83, 1).Value = "Cuba"
$ws.Cells.Item(83, 2).Value = 862
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 171
$ws.Cells.Item(83, 5).Value = 664
$ws.Cells.Item(83, 6).Value = 16
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 27

# Row 89: Letonia
$ws.Cells.Item(89, 1).Value = "Letonia"
$ws.Cells.Item(89, 2).Value = 682
$ws.Cells.Item(89, 3).Value = 7
$ws.Cells.Item(89, 4).Value = 57
$ws.Cells.Item(89, 5).Value = 620
$ws.Cells.Item(89, 6).Value = 5
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 5

# Row 90: Principado de Andorra
$ws.Cells.Item(90, 1).Value = "Principado de Andorra"
$ws.Cells.Item(90, 2).Value = 682
$ws.Cells.Item(90, 3).Value = 9
$ws.Cells.Item(90, 4).Value = 169
$ws.Cells.Item(90, 5).Value = 480
$ws.Cells.Item(90, 6).Value = 17
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 33

# Row 123: Islas Feroe
$ws.Cells.Item(123, 1).Value = "Islas Feroe"
$ws.Cells.Item(123, 2).Value = 184
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 171
$ws.Cells.Item(123, 5).Value = 13
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

# Row 190: Granada
$ws.Cells.Item(190, 1).Value = "Granada"
$ws.Cells.Item(190, 2).Value = 14
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 14
$ws.Cells.Item(190, 6).Value = 2
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 191: San Cristobal y Nieves
$ws.Cells.Item(191, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(191, 2).Value = 14
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 14
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0

# Row 195: Islas Malvinas
$ws.Cells.Item(195, 1).Value = "Islas Malvinas"
$ws.Cells.Item(195, 2).Value = 11
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 1
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

# Row 196: Islas Turcas y Caicos
$ws.Cells.Item(196, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(196, 2).Value = 11
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 10
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1

# Row 197: Montserrat
$ws.Cells.Item(197, 1).Value = "Montserrat"
$ws.Cells.Item(197, 2).Value = 11
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 1
$ws.Cells.Item(197, 5).Value = 10
$ws.Cells.Item(197, 6).Value = 1
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

# Row 198: Seychelles
$ws.Cells.Item(198, 1).Value = "Seychelles"
$ws.Cells.Item(198, 2).Value = 11
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 5
$ws.Cells.Item(198, 5).Value = 6
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0
